$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh -- update Price (D) and Volume(1h) (E) columns per row.
# D-column prices are text (t="inlineStr" in the source), and several look like
# plain numbers (e.g. "178.10", "1.00", "0.0000280"). A bare .Value assignment would
# let Excel auto-coerce those to a Double and silently drop the significant trailing
# zeros / split-decimal formatting, so each is written with a leading apostrophe to
# force text, then the cell style is reset to "Normal" so no stray quote-prefix / 
# text-format style sticks to the cell (keeps the original unstyled "General" cells).

$ws.Range("D2").Value = '''65.212.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.60%  '

$ws.Range("D3").Value = '''3.411.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.27%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''565.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.49%  '

$ws.Range("D6").Value = '''178.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.27%  '

$ws.Range("E7").Value = '  +4.31%  '

$ws.Range("D8").Value = '''3.401.71'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.29%  '

$ws.Range("D9").Value = '''0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("D10").Value = '''0.168'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +15.00%  '

$ws.Range("D11").Value = '''0.638'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.97%  '

$ws.Range("D12").Value = '''55.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.98%  '

$ws.Range("D13").Value = '''0.0000280'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.27%  '

$ws.Range("D14").Value = '''9.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.57%  '

$ws.Range("D15").Value = '''3.942.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("D16").Value = '''18.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.51%  '

$ws.Range("D17").Value = '''3.408.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.97%  '

$ws.Range("D18").Value = '''0.119'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.22%  '

$ws.Range("D19").Value = '''11.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.93%  '

$ws.Range("D20").Value = '''65.119.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.48%  '

$ws.Range("D21").Value = '''1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.49%  '

$ws.Range("D22").Value = '''469.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.75%  '

$ws.Range("D23").Value = '''5.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +18.83%  '

$ws.Range("D24").Value = '''4.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.67%  '

$ws.Range("D25").Value = '''86.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.40%  '

$ws.Range("D26").Value = '''13.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.15%  '

$ws.Range("D27").Value = '''10.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.32%  '

$ws.Range("E28").Value = '  +6.94%  '

$ws.Range("D29").Value = '''8.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.95%  '

$ws.Range("D30").Value = '''30.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.46%  '

$ws.Range("D31").Value = '''6.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.99%  '

$ws.Range("D32").Value = '''11.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.25%  '

$ws.Range("D33").Value = '''586.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.79%  '

$ws.Range("E34").Value = '  +5.40%  '

$ws.Range("D35").Value = '''60.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.84%  '

$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("E37").Value = '  -4.10%  '

$ws.Range("D38").Value = '''36.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.63%  '

$ws.Range("D39").Value = '''0.0₃0768'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.38%  '

$ws.Range("D40").Value = '''3.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").Value = '''0.377'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.18%  '

$ws.Range("D42").Value = '''3.119.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.96%  '

$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").Value = '''2.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.38%  '

$ws.Range("E45").Value = '  +3.89%  '

$ws.Range("D46").Value = '''0.0416'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.18%  '

$ws.Range("D47").Value = '''3.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.92%  '

$ws.Range("D48").Value = '''0.134'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.81%  '

$ws.Range("D49").Value = '''2.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("E50").Value = '  +7.18%  '

$ws.Range("D51").Value = '''137.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.74%  '
